$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 755.5
$ws.Range("I20").Value = 755.5
$ws.Range("K20").Value = 755.5
$ws.Range("M20").Value = -525.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 467.36365
$ws.Range("I33").Value = 264.2
$ws.Range("K33").Value = 264.2
$ws.Range("M33").Value = -35.19999999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 755.5
$ws.Range("I35").Value = 755.5
$ws.Range("K35").Value = 755.5
$ws.Range("M35").Value = -376.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 631
$ws.Range("I98").Value = 631
$ws.Range("K98").Value = 631
$ws.Range("M98").Value = 867

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 631
$ws.Range("I122").Value = 631
$ws.Range("K122").Value = 1893
$ws.Range("M122").Value = 557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 999
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1222.5714
$ws.Range("I137").Value = 1252.6666
$ws.Range("K137").Value = 3757.9998
$ws.Range("M137").Value = -1207.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2997.5
$ws.Range("I6").Value = 2997.5
$ws.Range("K6").Value = 2997.5
$ws.Range("M6").Value = -2824.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3914.8667
$ws.Range("I32").Value = 3496.3572
$ws.Range("K32").Value = 3496.3572
$ws.Range("M32").Value = -3209.3572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3455.6667
$ws.Range("I45").Value = 1240.4
$ws.Range("K45").Value = 1240.4
$ws.Range("M45").Value = -863.4000000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1359
$ws.Range("I61").Value = 1359
$ws.Range("K61").Value = 1359
$ws.Range("M61").Value = -1147

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1359
$ws.Range("I136").Value = 1359
$ws.Range("K136").Value = 4077
$ws.Range("M136").Value = -1527

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2923.4443
$ws.Range("I105").Value = 3195.25
$ws.Range("J105").Value = 749
$ws.Range("K105").Value = 3195.25
$ws.Range("L105").Value = 749
$ws.Range("M105").Value = -1448.25
$ws.Range("N105").Value = -4243

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 5835.5
$ws.Range("J106").Value = 5835.5
$ws.Range("L106").Value = 5835.5
$ws.Range("N106").Value = -8359.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 22280
$ws.Range("J52").Value = 22280
$ws.Range("L52").Value = 22280
$ws.Range("N52").Value = -22868

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 261
$ws.Range("I107").Value = 184.44444
$ws.Range("K107").Value = 184.44444
$ws.Range("M107").Value = 1735.55556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4833.1665
$ws.Range("I132").Value = 4833.1665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14499.4995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11969.4995
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 27120.2
$ws.Range("I22").Value = 36742.855
$ws.Range("J22").Value = 4667.3335
$ws.Range("K22").Value = 110228.565
$ws.Range("L22").Value = 14002.0005
$ws.Range("M22").Value = -110059.565
$ws.Range("N22").Value = -14340.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 27120.2
$ws.Range("I27").Value = 36742.855
$ws.Range("J27").Value = 4667.3335
$ws.Range("K27").Value = 110228.565
$ws.Range("L27").Value = 14002.0005
$ws.Range("M27").Value = -110126.565
$ws.Range("N27").Value = -14206.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 870.5714
$ws.Range("I113").Value = 300
$ws.Range("J113").Value = 965.6667
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 2897.0001
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -7237.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2888.1667
$ws.Range("J117").Value = 2965.8
$ws.Range("L117").Value = 8897.400000000001
$ws.Range("N117").Value = -15781.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1848
$ws.Range("I137").Value = 1323
$ws.Range("J137").Value = 3423
$ws.Range("K137").Value = 3969
$ws.Range("L137").Value = 10269
$ws.Range("M137").Value = 1131
$ws.Range("N137").Value = -20469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 8499
$ws.Range("I138").Value = 8499
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 25497
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -20357
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 10000
$ws.Range("I139").Value = 10000
$ws.Range("K139").Value = 30000
$ws.Range("M139").Value = -24860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 5324
$ws.Range("I107").Value = 377.14285
$ws.Range("J107").Value = 16866.666
$ws.Range("K107").Value = 377.14285
$ws.Range("L107").Value = 16866.666
$ws.Range("M107").Value = 1542.85715
$ws.Range("N107").Value = -20706.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4677.5884
$ws.Range("I126").Value = 3972.4167
$ws.Range("K126").Value = 11917.2501
$ws.Range("M126").Value = -9447.250100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 860.6
$ws.Range("I7").Value = 860.6
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 860.6
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -748.6
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 386.46667
$ws.Range("I55").Value = 133.33333
$ws.Range("J55").Value = 449.75
$ws.Range("K55").Value = 133.33333
$ws.Range("L55").Value = 449.75
$ws.Range("M55").Value = 39.66667000000001
$ws.Range("N55").Value = -795.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 860.6
$ws.Range("I126").Value = 860.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2581.8
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -111.8000000000002
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1680.1177
$ws.Range("I126").Value = 2016.4
$ws.Range("K126").Value = 6049.200000000001
$ws.Range("M126").Value = -3579.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 892.75
$ws.Range("I132").Value = 957.73334
$ws.Range("J132").Value = 697.8
$ws.Range("K132").Value = 2873.20002
$ws.Range("L132").Value = 2093.4
$ws.Range("M132").Value = -343.2000200000002
$ws.Range("N132").Value = -7153.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1224.3
$ws.Range("I136").Value = 1082.5555
$ws.Range("K136").Value = 3247.6665
$ws.Range("M136").Value = -697.6664999999998
